# Update view/favorite counts (column F) across the "展览", "演出" and
# "全部类型" sheets to match the regenerated GH-Pages data snapshot.
# ("本地生活" sheet is unaffected.)

$wb = $excel.ActiveWorkbook

# --- 展览 ("Exhibition") sheet ---------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7836
$ws1.Range("F5").Value = 7836
$ws1.Range("F7").Value = 120
$ws1.Range("F8").Value = 2124
$ws1.Range("F9").Value = 8609
$ws1.Range("F13").Value = 5750
$ws1.Range("F15").Value = 2733
$ws1.Range("F16").Value = 1197
$ws1.Range("F21").Value = 603
$ws1.Range("F22").Value = 48
$ws1.Range("F23").Value = 3869
$ws1.Range("F28").Value = 18
$ws1.Range("F29").Value = 5366
$ws1.Range("F32").Value = 288
$ws1.Range("F33").Value = 387
$ws1.Range("F34").Value = 158
$ws1.Range("F35").Value = 384
$ws1.Range("F36").Value = 2044
$ws1.Range("F37").Value = 1011
$ws1.Range("F39").Value = 1119
$ws1.Range("F40").Value = 4159
$ws1.Range("F43").Value = 34
$ws1.Range("F44").Value = 3501
$ws1.Range("F46").Value = 2335
$ws1.Range("F50").Value = 10

# --- 演出 ("Performance") sheet ---------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 68
$ws2.Range("F6").Value = 18
$ws2.Range("F10").Value = 7

# --- 全部类型 ("All types") sheet -------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7836
$ws4.Range("F5").Value = 7836
$ws4.Range("F7").Value = 120
$ws4.Range("F8").Value = 2124
$ws4.Range("F9").Value = 8609
$ws4.Range("F13").Value = 5750
$ws4.Range("F15").Value = 2733
$ws4.Range("F16").Value = 1197
$ws4.Range("F23").Value = 603
$ws4.Range("F25").Value = 48
$ws4.Range("F26").Value = 3869
$ws4.Range("F30").Value = 18
$ws4.Range("F31").Value = 5366
$ws4.Range("F33").Value = 387
$ws4.Range("F34").Value = 158
$ws4.Range("F35").Value = 384
$ws4.Range("F36").Value = 2044
$ws4.Range("F37").Value = 1011
$ws4.Range("F38").Value = 18
$ws4.Range("F40").Value = 1119
$ws4.Range("F42").Value = 4160
$ws4.Range("F45").Value = 34
$ws4.Range("F46").Value = 3501
$ws4.Range("F48").Value = 2335
